$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.315.96'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.002.23'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.60%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.87'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.69'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -6.59%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.522'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.003.21'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.60%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -6.47%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.73%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.64%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.59'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -7.13%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.503.63'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.26%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.290.37'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.004.56'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '454.50'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -5.50%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.679'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.31'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.16'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.28'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -6.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.28'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -5.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.02'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.16'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -4.34%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.89%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.85'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.16%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.96%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0792'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -6.02%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.95%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.12'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.05'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.92'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -11.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '409.81'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -6.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.276'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.29%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.772.93'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.16%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '38.16'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.85'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.83'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.13%  '
